# Auto-generated Excel COM-interop script applying scheduled-runner value updates
# to the Sheets workbook (Halicarnassus_Profits) per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 749.5
$ws.Range("I82").Value = 749.5
$ws.Range("K82").Value = 2248.5
$ws.Range("M82").Value = -1842.5

$ws.Range("H85").Value = 749.5
$ws.Range("I85").Value = 749.5
$ws.Range("K85").Value = 2248.5
$ws.Range("M85").Value = -844.5

$ws.Range("H87").Value = 77176.5
$ws.Range("J87").Value = 77176.5
$ws.Range("L87").Value = 77176.5
$ws.Range("N87").Value = -79672.5

$ws.Range("H90").Value = 77176.5
$ws.Range("J90").Value = 77176.5
$ws.Range("L90").Value = 231529.5
$ws.Range("N90").Value = -244009.5

$ws.Range("H112").Value = 1500
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -3392

$ws.Range("H116").Value = 8746.875
$ws.Range("I116").Value = 10229.167
$ws.Range("K116").Value = 10229.167
$ws.Range("M116").Value = -6787.166999999999

$ws.Range("H125").Value = 1863.7
$ws.Range("J125").Value = 3500
$ws.Range("L125").Value = 31500
$ws.Range("N125").Value = -36420

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3422.111
$ws.Range("I61").Value = 3385.5715
$ws.Range("J61").Value = 3550
$ws.Range("K61").Value = 3385.5715
$ws.Range("L61").Value = 3550
$ws.Range("M61").Value = -3173.5715
$ws.Range("N61").Value = -3974

$ws.Range("H88").Value = 1292.1111
$ws.Range("I88").Value = 851.5
$ws.Range("J88").Value = 1644.6
$ws.Range("K88").Value = 851.5
$ws.Range("L88").Value = 1644.6
$ws.Range("M88").Value = -445.5
$ws.Range("N88").Value = -2456.6

$ws.Range("H91").Value = 1292.1111
$ws.Range("I91").Value = 851.5
$ws.Range("J91").Value = 1644.6
$ws.Range("K91").Value = 851.5
$ws.Range("L91").Value = 1644.6
$ws.Range("M91").Value = 552.5
$ws.Range("N91").Value = -4452.6

$ws.Range("H136").Value = 3422.111
$ws.Range("I136").Value = 3385.5715
$ws.Range("J136").Value = 3550
$ws.Range("K136").Value = 10156.7145
$ws.Range("L136").Value = 10650
$ws.Range("M136").Value = -7606.7145
$ws.Range("N136").Value = -15750

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 171.88235
$ws.Range("I80").Value = 126
$ws.Range("J80").Value = 191
$ws.Range("K80").Value = 126
$ws.Range("L80").Value = 191
$ws.Range("M80").Value = 872
$ws.Range("N80").Value = -2187

$ws.Range("H83").Value = 171.88235
$ws.Range("I83").Value = 126
$ws.Range("J83").Value = 191
$ws.Range("K83").Value = 630
$ws.Range("L83").Value = 955
$ws.Range("M83").Value = 4362
$ws.Range("N83").Value = -10939

$ws.Range("H86").Value = 5431.125
$ws.Range("I86").Value = 2083
$ws.Range("K86").Value = 2083
$ws.Range("M86").Value = -960

$ws.Range("H89").Value = 5431.125
$ws.Range("I89").Value = 2083
$ws.Range("K89").Value = 10415
$ws.Range("M89").Value = -4799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8047.522
$ws.Range("J31").Value = 9361.941000000001
$ws.Range("L31").Value = 9361.941000000001
$ws.Range("N31").Value = -9951.941000000001

$ws.Range("H34").Value = 8047.522
$ws.Range("J34").Value = 9361.941000000001
$ws.Range("L34").Value = 9361.941000000001
$ws.Range("N34").Value = -9765.941000000001

$ws.Range("H58").Value = 3936.25
$ws.Range("I58").Value = 1323
$ws.Range("J58").Value = 7594.8
$ws.Range("K58").Value = 1323
$ws.Range("L58").Value = 7594.8
$ws.Range("M58").Value = -1120
$ws.Range("N58").Value = -8000.8

$ws.Range("H136").Value = 3936.25
$ws.Range("I136").Value = 1323
$ws.Range("J136").Value = 7594.8
$ws.Range("K136").Value = 3969
$ws.Range("L136").Value = 22784.4
$ws.Range("M136").Value = -1419
$ws.Range("N136").Value = -27884.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 524
$ws.Range("I70").Value = 524
$ws.Range("K70").Value = 1572
$ws.Range("M70").Value = -1257

$ws.Range("H73").Value = 524
$ws.Range("I73").Value = 524
$ws.Range("K73").Value = 1572
$ws.Range("M73").Value = -480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 29130.5
$ws.Range("J39").Value = 29130.5
$ws.Range("L39").Value = 29130.5
$ws.Range("N39").Value = -30194.5

$ws.Range("H70").Value = 9074.714
$ws.Range("I70").Value = 9128.75
$ws.Range("K70").Value = 9128.75
$ws.Range("M70").Value = -8858.75

$ws.Range("H73").Value = 9074.714
$ws.Range("I73").Value = 9128.75
$ws.Range("K73").Value = 9128.75
$ws.Range("M73").Value = -8192.75

$ws.Range("H80").Value = 1515.75
$ws.Range("I80").Value = 1019
$ws.Range("K80").Value = 1019
$ws.Range("M80").Value = -21

$ws.Range("H83").Value = 1515.75
$ws.Range("I83").Value = 1019
$ws.Range("K83").Value = 5095
$ws.Range("M83").Value = -103

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4714.1177
$ws.Range("I7").Value = 3196.1667
$ws.Range("K7").Value = 3196.1667
$ws.Range("M7").Value = -3084.1667

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = 0

$ws.Range("H47").Value = 5000
$ws.Range("I47").Value = 5000
$ws.Range("K47").Value = 5000
$ws.Range("M47").Value = -4510

$ws.Range("H52").Value = 5000
$ws.Range("I52").Value = 5000
$ws.Range("K52").Value = 5000
$ws.Range("M52").Value = -4767

$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51498

$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -157488

$ws.Range("H68").Value = 6349.5
$ws.Range("J68").Value = 10000
$ws.Range("L68").Value = 10000
$ws.Range("N68").Value = -11498

$ws.Range("H71").Value = 6349.5
$ws.Range("J71").Value = 10000
$ws.Range("L71").Value = 50000
$ws.Range("N71").Value = -57488

$ws.Range("H82").Value = 7598.1665
$ws.Range("I82").Value = 6794.5
$ws.Range("J82").Value = 8000
$ws.Range("K82").Value = 6794.5
$ws.Range("L82").Value = 8000
$ws.Range("M82").Value = -6433.5
$ws.Range("N82").Value = -8722

$ws.Range("H85").Value = 7598.1665
$ws.Range("I85").Value = 6794.5
$ws.Range("J85").Value = 8000
$ws.Range("K85").Value = 6794.5
$ws.Range("L85").Value = 8000
$ws.Range("M85").Value = -5546.5
$ws.Range("N85").Value = -10496

$ws.Range("H94").Value = 20000
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("N94").Value = -21352

$ws.Range("H98").Value = 63000
$ws.Range("J98").Value = 63000
$ws.Range("L98").Value = 63000
$ws.Range("N98").Value = -68990

$ws.Range("H100").Value = 4707.2173
$ws.Range("I100").Value = 1839.2
$ws.Range("J100").Value = 6913.385
$ws.Range("K100").Value = 1839.2
$ws.Range("L100").Value = 6913.385
$ws.Range("M100").Value = -1298.2
$ws.Range("N100").Value = -7995.385

$ws.Range("H126").Value = 4714.1177
$ws.Range("I126").Value = 3196.1667
$ws.Range("K126").Value = 9588.500100000001
$ws.Range("M126").Value = -7118.500100000001
